$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (BF) holds the game date that was mis-derived from the
# filename ("2-15-2013-14"). It should be the real ISO date (2014-02-15).
$dateCol = 58
$lastRow = 31

$oldValue = "2-15-2013-14"
$newValue = "2014-02-15"

$touched = $false
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Text -eq $oldValue) {
        # Assign via a string-literal formula so Excel's type inference keeps
        # this a plain text value instead of re-interpreting it as a date
        # serial (which a direct .Value/.Formula = "2014-02-15" would do).
        $cell.Formula = '="' + $newValue + '"'
        $touched = $true
    }
}

if ($touched) {
    $rng = $ws.Range($ws.Cells.Item(2, $dateCol), $ws.Cells.Item($lastRow, $dateCol))
    # Flatten the helper formulas down to literal string values so the saved
    # cells match plain text cells (no formula, no style/number-format churn).
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}
